# #78 adding changes for uploading match schedule in proper format
#
# Sheet1 previously held a header row (Home_Team, Away_Team, Tournament,
# Venue, matchTime, Description) followed by 3 rows of sample match data.
# The new format drops the header row + the Description column from
# Sheet1 (now just the 4 raw match rows in columns A:E) and adds a new
# Sheet2 that keeps the header row on top of the same 4 rows - i.e. the
# "proper format" for uploading the schedule.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- new match data (A:E => Home_Team, Away_Team, Tournament, Venue, matchTime)
$data = @(
    @("Mumbai Indians","Royal Challengers Bengaluru","IPL-20","Chinnaswamy Stadium","2020-08-27T10:47:20.868Z"),
    @("Chennai Super Kings","Royal Challengers Bengaluru","IPL-20","M. A. Chidambaram Stadium","2020-09-02T10:47:20.868Z"),
    @("Chennai Super Kings","Sun Risers Hyderabad","IPL-20","Feroz Shah Kotla","2020-09-05T10:47:20.868Z"),
    @("Chennai Super Kings","Mumbai Indians","IPL-20","Wankhede Stadium","2020-08-25T10:47:20.868Z")
)

# Wipe the old header/Description-column layout off Sheet1 and re-enter the
# 4 rows of match data only (no header row any more).
$ws1.Cells.Clear()
for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Add the new Sheet2 right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$header = @("Home_Team","Away_Team","Tournament","Venue","matchTime")
for ($c = 0; $c -lt $header.Length; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $header[$c]
}
for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}

# Carry the same column widths used on Sheet1 over to Sheet2.
$ws2.Columns.Item(1).ColumnWidth = 19.5867
$ws2.Columns.Item(2).ColumnWidth = 24.4187
$ws2.Columns.Item(3).ColumnWidth = 18.4187
$ws2.Columns.Item(5).ColumnWidth = 27.0867
$ws2.Columns.Item(6).ColumnWidth = 41.2507

# Restore the selection/active-sheet state: Sheet2 ends up with its whole
# grid selected, while Sheet1 stays the active tab with G13 selected.
$null = $ws2.Cells.Select()
$null = $ws1.Activate()
$null = $ws1.Range("G13").Select()
